$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Contest 5 (row 14, "KKR vs MI") results ---
$ws.Range("E14").Value2 = 40
$ws.Range("H14").Value2 = 0
$ws.Range("K14").Value2 = 20
$ws.Range("N14").Value2 = 100
$ws.Range("Q14").Value2 = 60
$ws.Range("T14").Value2 = 80

# --- Fill in match name for contest 11 (row 20) ---
$ws.Range("C20").Value2 = "DC vs SRH"

# --- Add two new contest rows (12 and 13) by inserting rows at 21:22,
#     matching the layout/style of row 20 ---
$ws.Rows("21:22").Insert()

$cols = @("A","B","C","D","E","G","H","J","K","M","N","P","Q","S","T")
foreach ($col in $cols) {
  $ws.Range("$col`20").Copy()
  $ws.Range("$col`21").PasteSpecial(-4122)
  $ws.Range("$col`22").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("A21").Value2 = 12
$ws.Range("B21").Value2 = 1
$ws.Range("C21").Value2 = "RR vs KKR"

$ws.Range("A22").Value2 = 13
$ws.Range("B22").Value2 = 1
$ws.Range("C22").Value2 = "KXI vs MI"

# --- Restore the selection that Excel shows in the saved file ---
$ws.Range("E29").Select()
